# Add six new weekly workout entries (rows 238-243) to Sheet1, mirroring a
# fresh Strava export appended below the existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column layout (row 1 header): A Participant | B Date | C Workout Type |
# D Total Duration | E Total Distance | F Total Elevation | G Zone 1 |
# H Zone 2 | I Zone 3 | J Zone 4 | K Zone 5 | L Workout Level | M Week

# Give the new date cells (column B) the same date number format already
# used by the column, by copying the formatted cell above them down first.
$ws.Range("B237").Copy($ws.Range("B238:B243"))

$rows = @(
    @{ Row=238; A="Matt";     B=45492; C="Ride";    D=76; E=21.24; F=1158; G=6;  H=33; I=26; J=8;  K=1;  L="Wily Hyena";   M=6 },
    @{ Row=239; A="Jeremiah"; B=45492; C="Workout"; D=47; E=0;     F=0;    G=44; H=2;  I=0;  J=0;  K=0;  L="Wily Hyena";   M=6 },
    @{ Row=240; A="Steven";   B=45492; C="Walk";    D=30; E=1.1;   F=36;   G=30; H=0;  I=0;  J=0;  K=0;  L="Brave Leopard"; M=6 },
    @{ Row=241; A="Steven";   B=45492; C="Walk";    D=34; E=1.71;  F=49;   G=34; H=0;  I=0;  J=0;  K=0;  L="Brave Leopard"; M=6 },
    @{ Row=242; A="Steven";   B=45493; C="Run";     D=30; E=2.85;  F=118;  G=1;  H=1;  I=19; J=8;  K=0;  L="Brave Leopard"; M=6 },
    @{ Row=243; A="Jeremiah"; B=45493; C="Workout"; D=61; E=0;     F=0;    G=61; H=0;  I=0;  J=0;  K=0;  L="Wily Hyena";   M=6 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
}

# Keep the frozen-pane view anchored near the bottom of the (now longer)
# table and move the active selection to the first blank row below it,
# matching how Excel leaves the cursor after a paste-append.
$ws.Range("A244").Select() | Out-Null
